$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table reports one column per year (D:R = 2007..2021). Add a new
# column S for year 2022 by mirroring column R (2021) - same formatting,
# new data.  xlPasteFormats = -4122.

# Row 2 (empty, bottom-bordered spacer cell) - copy format, no value
$ws.Range("R2").Copy()
$ws.Range("S2").PasteSpecial(-4122)

# Row 3 - year header 2022
$ws.Range("R3").Copy()
$ws.Range("S3").PasteSpecial(-4122)
$ws.Range("S3").Value = 2022

# Row 4 - population receiving pensions and disability benefits, persons
$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)
$ws.Range("S4").Value = 211650

# Row 5 - share of total population, percent
$ws.Range("R5").Copy()
$ws.Range("S5").PasteSpecial(-4122)
$ws.Range("S5").Value = 2.9794303052841493

$excel.CutCopyMode = 0

# Update the active selection to match the new last-used cell
$ws.Range("S2").Select()
